$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.510.56"
$ws.Range("E2").Value = "  -5.37%  "
$ws.Range("D3").Value = "2.470.51"
$ws.Range("E3").Value = "  -7.74%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.15"
$ws.Range("E5").Value = "  -2.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.85"
$ws.Range("E6").Value = "  -6.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  -4.00%  "
$ws.Range("D9").Value = "2.466.87"
$ws.Range("E9").Value = "  -7.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0995"
$ws.Range("E10").Value = "  -5.92%  "
$ws.Range("E11").Value = "  -2.62%  "
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.351"
$ws.Range("E13").Value = "  -4.43%  "
$ws.Range("D14").Value = "2.892.60"
$ws.Range("E14").Value = "  -8.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "24.02"
$ws.Range("E15").Value = "  -8.47%  "
$ws.Range("D16").Value = "59.410.37"
$ws.Range("E16").Value = "  -5.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000138"
$ws.Range("E17").Value = "  -6.03%  "
$ws.Range("D18").Value = "2.499.77"
$ws.Range("E18").Value = "  -6.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.14"
$ws.Range("E19").Value = "  -6.16%  "
$ws.Range("E20").Value = "  -5.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.71"
$ws.Range("E21").Value = "  -6.21%  "
$ws.Range("E22").Value = "  -3.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.74"
$ws.Range("E23").Value = "  -8.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.462"
$ws.Range("E24").Value = "  -8.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "60.67"
$ws.Range("E25").Value = "  -4.16%  "
$ws.Range("E26").Value = "  -4.27%  "
$ws.Range("E27").Value = "  -2.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.72"
$ws.Range("E28").Value = "  -5.76%  "
$ws.Range("E29").Value = "  -6.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.82"
$ws.Range("E30").Value = "  -6.11%  "
$ws.Range("D31").Value = "0.0₃0774"
$ws.Range("E31").Value = "  -9.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.26"
$ws.Range("E32").Value = "  -9.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.998"
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "155.27"
$ws.Range("E34").Value = "  -6.63%  "
$ws.Range("E35").Value = "  -6.03%  "
$ws.Range("E36").Value = "  -6.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.39"
$ws.Range("E37").Value = "  -5.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.74"
$ws.Range("E38").Value = "  -2.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.87"
$ws.Range("E39").Value = "  -6.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "316.09"
$ws.Range("E40").Value = "  -9.59%  "
$ws.Range("E41").Value = "  -4.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.839"
$ws.Range("E42").Value = "  -12.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.72"
$ws.Range("E43").Value = "  -7.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("E45").Value = "  -2.84%  "
$ws.Range("E46").Value = "  -5.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0939"
$ws.Range("E47").Value = "  -3.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0526"
$ws.Range("E48").Value = "  -6.76%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0229"
$ws.Range("E49").Value = "  -5.52%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.52"
$ws.Range("E50").Value = "  -8.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.91"
$ws.Range("E51").Value = "  -9.60%  "
